$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (Date serial, Method, ElapsedMs, wordCount, sentenceCount,
# posWordCount, negWordCount, posWordPercentage, negWordPercentage,
# positivePhraseCount, negativePhraseCount, posPhrasePercentage, negPhrasePercentage)
$rows = @(
    @(42600.792210648149, "Noun", 10884, 6475, 1129, 162, 60, 72, 27, 2, 1, 66, 33),
    @(42600.794548611113, "Noun", 13084, 7386, 1326, 175, 79, 68, 30, 2, 1, 66, 33),
    @(42600.830810185187, "Noun", 8717, 5450, 956, 110, 59, 64, 34, 2, 1, 66, 33),
    @(42600.879212962966, "Noun", 7086, 4367, 754, 72, 48, 59, 39, 1, 1, 50, 50)
)

$startRow = 3
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 1; $c -le $data.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}

# Widen column A slightly, as in the diff (13.85546875 -> 14.85546875)
$ws.Columns.Item(1).ColumnWidth = 14
